$d = $word.ActiveDocument

# 1. Rename the function definition: is_valid_password -> validate_password.
$d.Content.Find.Execute("def is_valid_password(password):", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "def validate_password(password):", 2)

# 2. The Test Lead's synced code block replaces the paragraphs that built the
#    valid/invalid password lists and looped asserts over them with seven
#    repeated "print(valid, msg)" lines (each followed by a blank paragraph).
#
#    Before (9 paragraphs, starting right after "def test_passwords():"):
#      "  valid_passwords = [...]"
#      "  invalid_passwords = [...]"
#      <blank>
#      "  for password in valid_passwords:"
#      "    assert valid, f\"...\""
#      <blank>
#      "  for password in invalid_passwords:"
#      "    assert not valid, f\"...\""
#      <blank>                                  <- left untouched, stays after the new block
#
#    After (14 paragraphs):
#      "  print(valid, msg)"   <blank>   (repeated 7 times)

$findStart = $d.Content.Find
$findStart.Text = "  valid_passwords = ["
$findStart.Execute() | Out-Null
$startPos = $findStart.Parent.Start

$findEnd = $d.Content.Find
$findEnd.Text = "    assert not valid, f`"Invalid password {password} passed validation`""
$findEnd.Execute() | Out-Null
$endPos = $findEnd.Parent.End

$blockRange = $d.Range($startPos, $endPos)

$printLine = "  print(valid, msg)"
$lines = @($printLine, $printLine, $printLine, $printLine, $printLine, $printLine, $printLine)
$blockRange.Text = [string]::Join("`r`r", $lines)
